# Append a new data row (row 19) to the merged stock/news sheet, mirroring
# the existing rows (date serial in A with the same date number format,
# numeric OHLCV/return/label values in B:H, and an empty headlines cell in I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

$ws.Cells.Item($row, 1).Value = 45859
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 2).Value = 6305.60009765625
$ws.Cells.Item($row, 3).Value = 6336.080078125
$ws.Cells.Item($row, 4).Value = 6303.7900390625
$ws.Cells.Item($row, 5).Value = 6304.740234375
$ws.Cells.Item($row, 6).Value = 5010840000
$ws.Cells.Item($row, 7).Value = 0.0013991348828683
$ws.Cells.Item($row, 8).Value = 1
$ws.Cells.Item($row, 9).Value = ""
